$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BVT")

# Insert two blank columns at L so the existing "Test Rail url" column
# (currently L) shifts right to N, freeing up L for the new field.
$ws.Columns("L:M").Insert()

# The hyperlink that lived on the old L2 still refers to the old address;
# re-anchor it on the relocated cell (now N2).
$ws.Hyperlinks.Delete(1)
$ws.Hyperlinks.Add($ws.Range("N2"), "https://surlatable.testrail.net/index.php?/cases/view/12080&group_by=cases:section_id&group_order=asc&display_deleted_cases=0&group_id=1961", "", "", "Test Rail url")

# Populate the newly freed "Quantity" column.
$ws.Range("L1").Value = "Quantity"
$ws.Range("L2").Value = 2

# Update the view: scroll so column I is the top-left visible column, and the
# active selection is the newly relocated hyperlink cell in N2.
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("N2").Select() | Out-Null
